$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet grows from 3 rows (1 header + 2 data) to 6 plain data rows
# (no header labels any more - row 1 becomes just another scraped
# record). We first blank out every cell that must end up empty using
# the "force text, then strip formatting" trick so the cell keeps an
# explicit empty-string value (and keeps the column/row inside the
# sheet's used range) instead of Excel garbage-collecting it back to a
# truly blank cell. Then we overwrite the cells that hold real values.
# ------------------------------------------------------------------

# Cells that must end up as explicit empty strings (columns A,B,C,F,G
# across all 6 rows, since those are never populated by the scraped
# D/E URL data).
$emptyCells = @(
  "A1","B1","C1","F1","G1",
  "C2","F2","G2",
  "C3","F3","G3",
  "A4","B4","C4","F4","G4",
  "A5","B5","C5","F5","G5",
  "A6","B6","C6","F6","G6"
)
foreach ($addr in $emptyCells) {
  $ws.Range($addr).Value = "'"
}

# Strip the quote-prefix/number-format noise that the trick above
# leaves behind, and also strip the old bold/centered header style
# that lived on row 1.
$ws.Range("A1:G6").ClearFormats()

# --- Row 1: scraped image data (was the header labels) ---
$ws.Range("D1").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty608/product/media/images/20221122/13/219567969/629009917/2/2_org_zoom.jpg"
$ws.Range("E1").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty608/product/media/images/20221122/13/219567969/629009917/2/2_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty607/product/media/images/20221122/13/219567969/629009917/3/3_org_zoom.jpg']"

# --- Row 2: new product replaces the shampoo row ---
$ws.Range("A2").Value = "S19 Max 4GB+64 GB Beyaz Cep Telefonu (Reeder Türkiye Garantili)"
$ws.Range("B2").Value = "3.049 TL"
$ws.Range("D2").Value = "https://cdn.dsmcdn.com/mnresize/1200/1800/ty844/product/media/images/20230425/13/331879740/844745791/1/1_org_zoom.jpg"
$ws.Range("E2").Value = "['https://cdn.dsmcdn.com/mnresize/1200/1800/ty844/product/media/images/20230425/13/331879740/844745791/1/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/seller-store/resources/hizli-satici.svg', 'https://cdn.dsmcdn.com/seller-store/resources/hizli-satici.svg', 'https://cdn.dsmcdn.com/ty844/product/media/images/20230425/13/331879740/844745791/1/1_org_zoom.jpg']"

# --- Row 3: new product replaces the shirt row ---
$ws.Range("A3").Value = "S19 Max 4GB+64 GB Mavi Cep Telefonu (Reeder Türkiye Garantili)"
$ws.Range("B3").Value = "2.999 TL"
$ws.Range("D3").Value = "https://cdn.dsmcdn.com/mnresize/1200/1800/ty847/product/media/images/20230426/11/332690742/918949896/1/1_org_zoom.jpg"
$ws.Range("E3").Value = "['https://cdn.dsmcdn.com/mnresize/1200/1800/ty847/product/media/images/20230426/11/332690742/918949896/1/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/ty847/product/media/images/20230426/11/332690742/918949896/1/1_org_zoom.jpg']"

# --- Row 4: brand new row ---
$ws.Range("D4").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty721/product/media/images/20230206/13/276465872/476930554/2/2_org_zoom.jpg"
$ws.Range("E4").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty721/product/media/images/20230206/13/276465872/476930554/2/2_org_zoom.jpg']"

# --- Row 5: brand new row ---
$ws.Range("D5").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty508/product/media/images/20220817/9/161294330/545780249/2/2_org_zoom.jpg"
$ws.Range("E5").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty508/product/media/images/20220817/9/161294330/545780249/2/2_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty507/product/media/images/20220817/9/161294330/545780249/3/3_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty507/product/media/images/20220817/9/161294330/545780249/4/4_org_zoom.jpg']"

# --- Row 6: brand new row ---
$ws.Range("D6").Value = "https://cdn.dsmcdn.com/mnresize/420/620/ty1150/product/media/images/prod/SPM/PIM/20240128/22/78a625f6-dc3d-3c43-9f55-7bf2c78824bd/1_org_zoom.jpg"
$ws.Range("E6").Value = "['https://cdn.dsmcdn.com/mnresize/420/620/ty1150/product/media/images/prod/SPM/PIM/20240128/22/78a625f6-dc3d-3c43-9f55-7bf2c78824bd/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1150/product/media/images/prod/SPM/PIM/20240128/22/93b8d939-c7f8-3225-b0c7-5b619c23490f/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/mnresize/420/620/ty1150/product/media/images/prod/SPM/PIM/20240128/22/a8abc4bf-a44a-3ab6-8da7-b75d0057aecc/1_org_zoom.jpg']"
